{"js": "const oldText = \"\u5e74\u30ad\u30e3\u30f3\u30da\u30fc\u30f3\u671f\u9593 \u5bfe\u8c61\uff1a\u306f\u304f\u3061\u3087\u3046\u5ea7 2022: 8\u670810\u301c19\u65e5\u30019\u67089\u301c18\u65e5\u300110\u67088\u301c17\u65e5\";\nconst newText = \" \uff1a2022\u5e74\u30ad\u30e3\u30f3\u30da\u30fc\u30f3\u671f\u9593 (\u5bfe\u8c61\uff1a\u306f\u304f\u3061\u3087\u3046\u5ea7)\uff1a\u30018\u670810\u301c19\u65e5\u30019\u67089\u301c18\u65e5\u300110\u67088\u301c17\u65e5\";\n\nconst results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"\u5e74\u30ad\u30e3\u30f3\u30da\u30fc\u30f3\u671f\u9593 \u5bfe\u8c61\uff1a\u306f\u304f\u3061\u3087\u3046\u5ea7 2022: 8\u670810\u301c19\u65e5\u30019\u67089\u301c18\u65e5\u300110\u67088\u301c17\u65e5\"\n$newText = \" \uff1a2022\u5e74\u30ad\u30e3\u30f3\u30da\u30fc\u30f3\u671f\u9593 (\u5bfe\u8c61\uff1a\u306f\u304f\u3061\u3087\u3046\u5ea7)\uff1a\u30018\u670810\u301c19\u65e5\u30019\u67089\u301c18\u65e5\u300110\u67088\u301c17\u65e5\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Replacement.ClearFormatting()\n$range.Find.Execute(\n    $oldText,        # FindText\n    $false,          # MatchCase\n    $false,          # MatchWholeWord\n    $false,          # MatchWildcards\n    $false,          # MatchSoundsLike\n    $false,          # MatchAllWordForms\n    $true,           # Forward\n    $wdFindContinue, # Wrap\n    $false,          # Format\n    $newText,        # ReplaceWith\n    $wdReplaceAll    # Replace\n) | Out-Null\n"}
